# Updated cryptos list on Thu Jul  4 18:46:29 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.169.28"
$ws.Range("E2").Value = "  -3.68%  "

# Row 3
$ws.Range("D3").Value = "3.134.33"
$ws.Range("E3").Value = "  -5.17%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "'523.65"
$ws.Range("E5").Value = "  -6.21%  "

# Row 6
$ws.Range("D6").Value = "'134.32"
$ws.Range("E6").Value = "  -5.55%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "3.134.18"
$ws.Range("E8").Value = "  -5.14%  "

# Row 9
$ws.Range("E9").Value = "  -6.44%  "

# Row 10
$ws.Range("D10").Value = "'7.19"
$ws.Range("E10").Value = "  -8.38%  "

# Row 11
$ws.Range("E11").Value = "  -9.25%  "

# Row 12
$ws.Range("E12").Value = "  -7.58%  "

# Row 13
$ws.Range("D13").Value = "3.674.83"
$ws.Range("E13").Value = "  -5.11%  "

# Row 14
$ws.Range("E14").Value = "  -1.84%  "

# Row 15
$ws.Range("D15").Value = "'25.39"
$ws.Range("E15").Value = "  -5.65%  "

# Row 16
$ws.Range("D16").Value = "3.136.86"
$ws.Range("E16").Value = "  -5.14%  "

# Row 17
$ws.Range("D17").Value = "58.170.22"
$ws.Range("E17").Value = "  -3.69%  "

# Row 18
$ws.Range("D18").Value = "'0.0000151"
$ws.Range("E18").Value = "  -8.62%  "

# Row 19
$ws.Range("D19").Value = "'5.77"
$ws.Range("E19").Value = "  -5.55%  "

# Row 20
$ws.Range("D20").Value = "'13.01"
$ws.Range("E20").Value = "  -6.90%  "

# Row 21
$ws.Range("D21").Value = "'7.89"
$ws.Range("E21").Value = "  -8.77%  "

# Row 22
$ws.Range("D22").Value = "'342.18"
$ws.Range("E22").Value = "  -8.84%  "

# Row 23
$ws.Range("E23").Value = "  -0.15%  "

# Row 24
$ws.Range("D24").Value = "'67.69"
$ws.Range("E24").Value = "  -9.00%  "

# Row 25
$ws.Range("D25").Value = "'0.504"
$ws.Range("E25").Value = "  -5.82%  "

# Row 26
$ws.Range("D26").Value = "3.263.85"
$ws.Range("E26").Value = "  -5.04%  "

# Row 27
$ws.Range("D27").Value = "'0.166"
$ws.Range("E27").Value = "  -3.34%  "

# Row 28
$ws.Range("D28").Value = "0.0{0}0952" -f [char]0x2083
$ws.Range("E28").Value = "  -7.16%  "

# Row 29
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.05%  "

# Row 30
$ws.Range("D30").Value = "'6.78"
$ws.Range("E30").Value = "  -5.50%  "

# Row 31
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("E32").Value = "  -8.96%  "

# Row 33
$ws.Range("E33").Value = "  -10.23%  "

# Row 34
$ws.Range("D34").Value = "'21.36"
$ws.Range("E34").Value = "  -5.64%  "

# Row 35
$ws.Range("E35").Value = "  -3.00%  "

# Row 36
$ws.Range("D36").Value = "'4.77"
$ws.Range("E36").Value = "  -6.79%  "

# Row 37
$ws.Range("D37").Value = "'156.53"
$ws.Range("E37").Value = "  -6.22%  "

# Row 38
$ws.Range("D38").Value = "'6.20"
$ws.Range("E38").Value = "  -7.19%  "

# Row 39
$ws.Range("E39").Value = "  -11.18%  "

# Row 40
$ws.Range("D40").Value = "'0.0684"
$ws.Range("E40").Value = "  -6.44%  "

# Row 41
$ws.Range("D41").Value = "3.168.75"
$ws.Range("E41").Value = "  -5.07%  "

# Row 42
$ws.Range("D42").Value = "'40.33"
$ws.Range("E42").Value = "  -3.83%  "

# Row 43
$ws.Range("D43").Value = "'24.00"
$ws.Range("E43").Value = "  -9.91%  "

# Row 44
$ws.Range("D44").Value = "'0.692"
$ws.Range("E44").Value = "  -7.78%  "

# Row 45
$ws.Range("E45").Value = "  -3.04%  "

# Row 46
$ws.Range("D46").Value = "'3.88"
$ws.Range("E46").Value = "  -6.08%  "

# Row 47
$ws.Range("E47").Value = "  +0.00%  "

# Row 48
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.259.17"
$ws.Range("E48").Value = "  -3.93%  "

# Row 49
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.43"
$ws.Range("E49").Value = "  -9.34%  "

# Row 50
$ws.Range("D50").Value = "'6.18"
$ws.Range("E50").Value = "  -3.21%  "

# Row 51
$ws.Range("D51").Value = "'20.55"
$ws.Range("E51").Value = "  -3.72%  "

